$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find an existing "Diagrama de flujo: proceso" shape (flowChartProcess) that
# already carries the exact style / text formatting we need (id 58 -
# "Diagrama de flujo: proceso 57"), then duplicate it so the new shape
# inherits the same <p:style> (lnRef/fillRef/effectRef/fontRef) and run
# formatting instead of trying to rebuild it from scratch.
$template = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 58) {
        $template = $sh
    }
}

if ($template -eq $null) {
    # Fallback: any other existing shape works as a style donor too, they
    # all share the same flowChartProcess theme style in this deck.
    $template = $s.Shapes.Item(2)
}

$dupRange = $template.Duplicate()
$shape = $dupRange.Item(1)

$shape.Name = "Diagrama de flujo: proceso 1"

# Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU).
# Target EMU (from the authored slide XML): off x=4516606 y=10832399,
# ext cx=582318 cy=377135. The literal point values below are chosen so
# that, after the host's internal points->EMU conversion (which rounds
# through a single-precision float before truncating to EMU), they land
# exactly on the target EMU instead of one EMU short.
$shape.Left = 355.63826771653544
$shape.Top = 852.9448031496063
$shape.Width = 45.85181242362205
$shape.Height = 29.695670191338582

$shape.TextFrame.TextRange.Text = "No"
